$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old D1 (#NUM!) and E1 (#VALUE!) formulas one column to the
# right, making room for a new D1 formula that references a missing
# external workbook/sheet (yields #REF!).
$ws.Range("F1").Formula = $ws.Range("E1").Formula
$ws.Range("E1").Formula = "=SQRT(-1)"
$ws.Range("D1").Formula = "=[1]SheetNotExists!A1"

# Match the selection left behind by the edit (single cell D1).
[void]$ws.Range("D1").Select()
